# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The only observable change is the "K" column (column G) values for rows 2-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 1
    20 = 3
    21 = 1
    22 = 2
    23 = 0
    24 = 0
    25 = 1
    26 = 1
    27 = 1
    28 = 0
    29 = 1
    30 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
